# Apply the edits described by the diff:
# 1. Rename the worksheet/tab from "GossA-HW40.xpc" to "GossA"
# 2. Append a new row 16 with averaged-intensity data for the
#    "HexGrid-60degTilt5degRes" Gaussian-Quadrature scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "GossA"

# 2. Add new row 16 data (mirrors the layout of row 15)
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 1.002854316891267
$ws.Cells.Item(16, 4).Value = 0.9872471667522974
$ws.Cells.Item(16, 5).Value = 1.001031114751944
$ws.Cells.Item(16, 6).Value = 0.9959984606489831
$ws.Cells.Item(16, 7).Value = 1.002854316891267
$ws.Cells.Item(16, 8).Value = 0.9872471667522974
$ws.Cells.Item(16, 9).Value = 0.9992449151762937
$ws.Cells.Item(16, 10).Value = 0.9964413401068091
$ws.Cells.Item(16, 11).Value = 1.000508245343851
$ws.Cells.Item(16, 12).Value = 0.9935001636362208
$ws.Cells.Item(16, 13).Value = 1.002854316891267
$ws.Cells.Item(16, 14).Value = 0.9941391407521205
$ws.Cells.Item(16, 15).Value = 0.9967827647611227
$ws.Cells.Item(16, 16).Value = 0.9971032154134581

# Copy the formatting of A15 (bordered/bold/centered "index" style) onto A16,
# then set its value -- mirrors how the preceding index column cells look.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Value = 14

$excel.CutCopyMode = $false
